# FEAT: better handling of % and %% as words.
# Adds a new lexer-FSM state S_PERCENT (row for the "%" state) to the
# `transitions` sheet, wires a couple of existing states to transition into
# it on a "%" character, and fixes a couple of cells that were previously
# marked as T_ERROR but should legally continue as T_WORD.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write a cell's text and stamp it with the formatting of another
# cell that already carries the desired style (mirrors "paste values then
# copy format" which is how these FSM tables were actually edited).
function Set-CellWithStyle {
    param($sheet, [string]$CellRef, [string]$Text, [string]$StyleFromRef)
    $sheet.Range($CellRef).Value = $Text
    $sheet.Range($StyleFromRef).Copy() | Out-Null
    $sheet.Range($CellRef).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $excel.CutCopyMode = 0
}

# --- Row 8 (S_FILE_1ST state): O8 was wrongly T_WORD, should be S_PERCENT ---
Set-CellWithStyle $ws 'O8' 'S_PERCENT' 'N8'

# --- Row 13 (S_HDPER_ST state): several T_ERROR cells are actually legal T_WORD ---
foreach ($col in @('B','C','F','G','H','I','K','L','X','AD','AN')) {
    Set-CellWithStyle $ws "${col}13" 'T_WORD' 'F2'
}

# --- Row 51 (S_WORD_1ST state): on "%" (col AB) now transition into S_PERCENT ---
Set-CellWithStyle $ws 'AB51' 'S_PERCENT' 'M53'

# --- Insert the new S_PERCENT state as row 54 (pushes old rows 54-61 to 55-62) ---
$ws.Rows.Item(54).Insert()

# Row label + per-character-class transitions for the new S_PERCENT state.
# Default target is T_ERROR, except: word-continuation classes go to
# T_WORD, and another "%" (col AB) loops back into S_PERCENT itself.
Set-CellWithStyle $ws 'A54' 'S_PERCENT' 'A3'

$wordCols = @('B','C','F','G','H','I','J','K','X','AD','AN')
foreach ($col in $wordCols) {
    Set-CellWithStyle $ws "${col}54" 'T_WORD' 'F2'
}

$errorCols = @('D','E','L','M','N','O','P','Q','R','S','T','U','V','W','Y','Z','AA','AC','AE','AF','AG','AH','AI','AJ','AK','AL','AM')
foreach ($col in $errorCols) {
    Set-CellWithStyle $ws "${col}54" 'T_ERROR' 'K2'
}

Set-CellWithStyle $ws 'AB54' 'S_PERCENT' 'M53'

# --- Recompute row heights for every row we touched (Excel drops the custom
#     height once a row's content/format changes and it no longer needs it) ---
foreach ($r in @(8, 13, 51, 53, 54)) {
    $ws.Rows.Item($r).AutoFit()
}

# --- Leave the selection where the edit actually happened ---
$ws.Range('O8').Select() | Out-Null
